# "Send mails working now"
# - Replace the Faculty Mentor ("Akshit Garg") entries for rows 2, 6, 8, 11
#   with the new mentor "be18103032 Gaganpreet Singh Khurana" (a brand new
#   shared string gets created for this).
# - Move the active selection to G14 (from L10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newMentor = "be18103032 Gaganpreet Singh Khurana"

$ws.Range("L2").Value = $newMentor
$ws.Range("L6").Value = $newMentor
$ws.Range("L8").Value = $newMentor
$ws.Range("L11").Value = $newMentor

[void]$ws.Range("G14").Select()
